$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.057.27'
$ws.Range('E2').Value = '''  -5.66%  '
$ws.Range('D3').Value = '''3.570.13'
$ws.Range('E3').Value = '''  -0.64%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '''  +0.07%  '
$ws.Range('D5').Value = '''400.13'
$ws.Range('E5').Value = '''  -3.72%  '
$ws.Range('D6').Value = '''131.02'
$ws.Range('E6').Value = '''  +0.54%  '
$ws.Range('D7').Value = '''3.719.86'
$ws.Range('E7').Value = '''  +3.72%  '
$ws.Range('D8').Value = '''0.607'
$ws.Range('E8').Value = '''  -6.40%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '''  -0.09%  '
$ws.Range('D10').Value = '''0.708'
$ws.Range('E10').Value = '''  -8.79%  '
$ws.Range('D11').Value = '''0.156'
$ws.Range('E11').Value = '''  -10.88%  '
$ws.Range('D12').Value = '''0.0000308'
$ws.Range('E12').Value = '''  -8.71%  '
$ws.Range('D13').Value = '''40.79'
$ws.Range('E13').Value = '''  -4.10%  '
$ws.Range('D14').Value = '''9.68'
$ws.Range('E14').Value = '''  -2.26%  '
$ws.Range('D15').Value = '''4.168.54'
$ws.Range('E15').Value = '''  +0.03%  '
$ws.Range('E16').Value = '''  -1.43%  '
$ws.Range('D17').Value = '''3.584.09'
$ws.Range('E17').Value = '''  -1.27%  '
$ws.Range('D18').Value = '''19.56'
$ws.Range('E18').Value = '''  -3.90%  '
$ws.Range('D19').Value = '''13.08'
$ws.Range('E19').Value = '''  +5.59%  '
$ws.Range('D20').Value = '''1.06'
$ws.Range('E20').Value = '''  -7.21%  '
$ws.Range('D21').Value = '''63.236.25'
$ws.Range('E21').Value = '''  -5.19%  '
$ws.Range('D22').Value = '''411.88'
$ws.Range('E22').Value = '''  -8.41%  '
$ws.Range('D23').Value = '''14.77'
$ws.Range('E23').Value = '''  +12.80%  '
$ws.Range('D24').Value = '''83.87'
$ws.Range('E24').Value = '''  -5.81%  '
$ws.Range('D25').Value = '''2.93'
$ws.Range('E25').Value = '''  -7.14%  '
$ws.Range('D26').Value = '''34.85'
$ws.Range('E26').Value = '''  -1.27%  '
$ws.Range('D27').Value = '''3.13'
$ws.Range('E27').Value = '''  -6.19%  '
$ws.Range('D28').Value = '''9.20'
$ws.Range('E28').Value = '''  -7.57%  '
$ws.Range('D29').Value = '''5.08'
$ws.Range('E29').Value = '''  +4.87%  '
$ws.Range('D30').Value = '''12.44'
$ws.Range('E30').Value = '''  +0.48%  '
$ws.Range('D31').Value = '''2.68'
$ws.Range('E31').Value = '''  -2.80%  '
$ws.Range('D32').Value = '''0.113'
$ws.Range('E32').Value = '''  -3.39%  '
$ws.Range('D33').Value = '''6.80'
$ws.Range('E33').Value = '''  -7.86%  '
$ws.Range('D34').Value = '''0.157'
$ws.Range('E34').Value = '''  -2.72%  '
$ws.Range('D35').Value = '''39.88'
$ws.Range('E35').Value = '''  -1.51%  '
$ws.Range('D36').Value = '''0.998'
$ws.Range('E36').Value = '''  -0.16%  '
$ws.Range('D37').Value = '''55.22'
$ws.Range('D38').Value = '''0.0454'
$ws.Range('E38').Value = '''  -8.16%  '
$ws.Range('B39').Value = '''EnergySwap'
$ws.Range('C39').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '''28.43'
$ws.Range('E39').Value = '''  +32.53%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '''  +0.49%  '
$ws.Range('B41').Value = '''ThetaToken'
$ws.Range('C41').Value = '''https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = '''2.82'
$ws.Range('E41').Value = '''  +22.33%  '
$ws.Range('B42').Value = '''Stellar'
$ws.Range('C42').Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '''0.138'
$ws.Range('E42').Value = '''  -6.29%  '
$ws.Range('B43').Value = '''ApeXProtocol'
$ws.Range('C43').Value = '''https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = '''3.12'
$ws.Range('E43').Value = '''  +21.91%  '
$ws.Range('B44').Value = '''Monero'
$ws.Range('C44').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''143.54'
$ws.Range('E44').Value = '''  -4.01%  '
$ws.Range('B45').Value = '''NEARProtocol'
$ws.Range('C45').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''4.30'
$ws.Range('E45').Value = '''  -0.59%  '
$ws.Range('D46').Value = '''3.23'
$ws.Range('E46').Value = '''  -1.93%  '
$ws.Range('B47').Value = '''PEPE'
$ws.Range('C47').Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '''0.0₃0623'
$ws.Range('E47').Value = '''  -13.86%  '
$ws.Range('D48').Value = '''2.02'
$ws.Range('E48').Value = '''  +2.11%  '
$ws.Range('D49').Value = '''2.76'
$ws.Range('E49').Value = '''  -8.29%  '
$ws.Range('D50').Value = '''2.51'
$ws.Range('E50').Value = '''  -8.46%  '
$ws.Range('D51').Value = '''0.284'
$ws.Range('E51').Value = '''  -9.73%  '
